$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing QA rows (2-4) and PREPROD rows (5-8) with new claim numbers
# Leading apostrophe forces Excel to treat the numeric-looking string as text
# (preserving leading zeros / trailing spaces) while keeping the quotePrefix style.
$ws.Range("B2").Value = "'1120194100440   "
$ws.Range("B3").Value = "'1220194200683"
$ws.Range("B4").Value = "'0420194406812"
$ws.Range("B5").Value = "'0420172008629    "
$ws.Range("B6").Value = "'0420172008636"
$ws.Range("B7").Value = "'0420172008630"
$ws.Range("B8").Value = "'0420172008637"

# Add new rows 9, 10, 11 for the two additional "smart folders" (PREPROD)
$ws.Range("A9").Value = "PREPROD"
$ws.Range("B9").Value = "'1220170301441"
$ws.Range("B9").QuotePrefix = $true

$ws.Range("A10").Value = "PREPROD"
$ws.Range("B10").Value = "'1120170200942"
$ws.Range("B10").QuotePrefix = $true

$ws.Range("A11").Value = "PREPROD"
$ws.Range("B11").Value = "'1220170301442 "
$ws.Range("B11").QuotePrefix = $true

# Update selection to match the new active cell
$ws.Range("A11").Select()
